# Inserts one new data row for "Cilantro" at sheet row 87 (pushing the
# existing rows 87-211 down to 88-212), matching the target diff which
# grows the sheet from A1:R211 to A1:R212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 87, shifting everything below
# down by one (Excel's normal "insert row" behaviour).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Cells.Item(87, 1).Value = 8
$ws.Cells.Item(87, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(87, 3).Value = 'Coquimbo'
$ws.Cells.Item(87, 4).Value = 44994
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100112040
$ws.Cells.Item(87, 7).Value = 'Cilantro'
$ws.Cells.Item(87, 8).Value = 'Sin especificar'
$ws.Cells.Item(87, 9).Value = 'Primera'
$ws.Cells.Item(87, 10).Value = 2100
$ws.Cells.Item(87, 11).Value = 2000
$ws.Cells.Item(87, 12).Value = 2500
$ws.Cells.Item(87, 13).Value = 2250
$ws.Cells.Item(87, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(87, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(87, 16).Value = 1500
$ws.Cells.Item(87, 17).Value = 1.5
$ws.Cells.Item(87, 18).Value = 'Hortaliza'

# Make sure the date column keeps the same date/time number format the
# rest of column D uses (copy format down from the row above).
$ws.Cells.Item(86, 4).Copy()
$ws.Cells.Item(87, 4).PasteSpecial(-4122)
$ws.Cells.Item(87, 4).Value = 44994
